$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The metaParamTex column (T) for the two "Stylized-Normal" rows and the two
# "Log-Normal" rows (rows 5-8) should use the LaTeX symbol for the mean,
# "\mu", instead of the generic "\beta" placeholder that was there before.
$ws.Range("T5:T8").Value = "\mu"

# Reflect the refactor in the sheet's active selection: focus moves from the
# previously-selected Poisson/Exponential slider rows (T9:T14) to the rows
# that were just updated (T5:T8), with T5 as the active cell.
[void]$ws.Range("T5:T8").Select()
